# Update to-do for each deliverable
# ------------------------------------------------------------
# Adds a "Need to do review" note (wrapped) to several rows of
# the Checklist sheet's "Notes" column (G), marks a few more rows
# as "Ready for printing", enlarges column G, tightens row
# heights, and adds a conditional-format rule that highlights
# "Ready" text in column G.
# ------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Checklist"

$reviewNote = "Need to do review :`n+ review form (if there's changes)`n+ audit meeting minutes`n+ Update recordlogs.xls with the new file location"
$readyText  = "Ready for printing"

# ---- Column G (Notes) updates -------------------------------------------
# Rows that get the multi-line "Need to do review" note (wrapped text)
$reviewRows = @(2, 6, 7, 8, 9)
foreach ($r in $reviewRows) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value2 = $reviewNote
    $cell.WrapText = $true
}

# Row that gets the "Ready for printing" note (re-uses existing text/style)
$ws.Cells.Item(3, 7).Value2 = $readyText

# ---- Row heights ----------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 14.25
$ws.Rows.Item(2).RowHeight = 63.75
$ws.Rows.Item(3).RowHeight = 12.75
$ws.Rows.Item(4).RowHeight = 12.75
$ws.Rows.Item(5).RowHeight = 12.75
$ws.Rows.Item(6).RowHeight = 63.75
$ws.Rows.Item(7).RowHeight = 63.75
$ws.Rows.Item(8).RowHeight = 63.75
$ws.Rows.Item(9).RowHeight = 63.75
$ws.Rows.Item(10).RowHeight = 12.75
$ws.Rows.Item(12).RowHeight = 12.75
$ws.Rows.Item(13).RowHeight = 12.75
$ws.Rows.Item(14).RowHeight = 12.75
$ws.Rows.Item(15).RowHeight = 25.5
$ws.Rows.Item(16).RowHeight = 12.75
$ws.Rows.Item(17).RowHeight = 12.75

# ---- Column G width ---------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 35.3

# ---- Conditional formatting: highlight "Ready" in column G ------------
$rng = $ws.Range("G1:G1048576")
$fc = $rng.FormatConditions.Add(9, 0, "Ready")
$fc.Text = "Ready"
$fc.Formula1 = 'NOT(ISERROR(SEARCH("Ready",G1)))'
$fc.Font.Color = 24832
$fc.Interior.Color = 13561798
$fc.SetFirstPriority()

# ---- Selection / scroll position ---------------------------------------
$ws.Activate()
$ws.Range("C11").Select()
